# Update the cryptos price/volume list with the latest values.
# Values that look like plain numbers (single '.') are prefixed with an
# apostrophe so Excel stores them as text (matching the workbook's existing
# "D.D" style text formatting) rather than auto-converting them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = '67.460.03'
$ws.Range("E2").Value = '  +4.06%  '
$ws.Range("D3").Value = '3.255.51'
$ws.Range("E3").Value = '  +3.86%  '
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").Value = "'578.05"
$ws.Range("E5").Value = '  +2.16%  '
$ws.Range("D6").Value = "'181.80"
$ws.Range("E6").Value = '  +8.33%  '
$ws.Range("E7").Value = '  -0.05%  '
$ws.Range("D8").Value = "'0.598"
$ws.Range("E8").Value = '  -1.66%  '
$ws.Range("D9").Value = '3.252.01'
$ws.Range("E9").Value = '  +3.69%  '
$ws.Range("E10").Value = '  +6.71%  '
$ws.Range("E11").Value = '  +3.20%  '
$ws.Range("E12").Value = '  +6.74%  '
$ws.Range("D13").Value = '3.816.02'
$ws.Range("E13").Value = '  +3.67%  '
$ws.Range("E14").Value = '  +1.06%  '
$ws.Range("D15").Value = "'28.47"
$ws.Range("E15").Value = '  +6.45%  '
$ws.Range("D16").Value = '67.421.31'
$ws.Range("E16").Value = '  +4.01%  '
$ws.Range("D17").Value = "'0.0000168"
$ws.Range("E17").Value = '  +4.46%  '
$ws.Range("D18").Value = '3.252.17'
$ws.Range("E18").Value = '  +3.39%  '
$ws.Range("E19").Value = '  +2.99%  '
$ws.Range("D20").Value = "'13.56"
$ws.Range("E20").Value = '  +6.94%  '
$ws.Range("D21").Value = "'375.83"
$ws.Range("E21").Value = '  +5.92%  '
$ws.Range("D22").Value = "'7.62"
$ws.Range("E22").Value = '  +6.22%  '
$ws.Range("E23").Value = '  -0.05%  '
$ws.Range("D24").Value = "'71.11"
$ws.Range("E24").Value = '  +3.21%  '
$ws.Range("E25").Value = '  +4.13%  '
$ws.Range("E26").Value = '  +4.48%  '
$ws.Range("D27").Value = "'9.55"
$ws.Range("E27").Value = '  -0.23%  '
$ws.Range("E28").Value = '  +2.50%  '
$ws.Range("E29").Value = '  -0.22%  '
$ws.Range("E30").Value = '  +9.00%  '
$ws.Range("E31").Value = '  +3.57%  '
$ws.Range("D32").Value = "'22.68"
$ws.Range("E32").Value = '  +4.25%  '
$ws.Range("E33").Value = '  +0.07%  '
$ws.Range("E34").Value = '  +6.79%  '
$ws.Range("E35").Value = '  +5.83%  '
$ws.Range("D36").Value = "'163.48"
$ws.Range("E36").Value = '  +3.52%  '
$ws.Range("E37").Value = '  +5.33%  '
$ws.Range("E38").Value = '  +2.22%  '
$ws.Range("E39").Value = '  +5.80%  '
$ws.Range("D40").Value = "'6.85"
$ws.Range("E40").Value = '  +13.28%  '
$ws.Range("E41").Value = '  +12.92%  '
$ws.Range("D42").Value = "'26.62"
$ws.Range("E42").Value = '  +3.38%  '
$ws.Range("E43").Value = '  +8.35%  '
$ws.Range("D44").Value = "'358.10"
$ws.Range("E44").Value = '  +13.02%  '
$ws.Range("D45").Value = '2.705.72'
$ws.Range("E45").Value = '  +2.12%  '
$ws.Range("D46").Value = "'25.41"
$ws.Range("E46").Value = '  +6.91%  '
$ws.Range("D47").Value = "'40.81"
$ws.Range("E47").Value = '  +3.65%  '
$ws.Range("D48").Value = "'0.0680"
$ws.Range("E48").Value = '  +4.91%  '
$ws.Range("E49").Value = '  +3.70%  '
$ws.Range("D50").Value = "'0.998"
$ws.Range("E50").Value = '  +7.45%  '
$ws.Range("E51").Value = '  -0.19%  '
